# This script reproduces the crypto price/volume refresh described by the commit
# "Updated cryptos list ... with GitHub Actions": columns D (Price) and E
# (Volume(1h)) are refreshed for every coin row (rows 2-51).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.007.13'
$ws.Range("E2").Value = '  +2.20%  '
$ws.Range("D3").Value = '1.906.53'
$ws.Range("E3").Value = '  +2.51%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -0.78%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '317.55'
$ws.Range("E5").Value = '  +2.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.003'
$ws.Range("E6").Value = '  -0.72%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4821'
$ws.Range("E7").Value = '  +1.17%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07360'
$ws.Range("E9").Value = '  +0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9326'
$ws.Range("E10").Value = '  +0.21%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.83'
$ws.Range("E11").Value = '  +0.15%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07749'
$ws.Range("E12").Value = '  -0.51%  '
$ws.Range("D13").Value = '1.920.98'
$ws.Range("E13").Value = '  +3.44%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.494'
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.636'
$ws.Range("E15").Value = '  +1.47%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.86'
$ws.Range("E16").Value = '  +2.08%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.005'
$ws.Range("E17").Value = '  -0.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008889'
$ws.Range("E18").Value = '  +0.90%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.003'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").Value = '28.036.39'
$ws.Range("E20").Value = '  +2.38%  '
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.146'
$ws.Range("E22").Value = '  +1.01%  '
$ws.Range("D23").Value = '2.139.20'
$ws.Range("E23").Value = '  +2.36%  '
$ws.Range("E24").Value = '  +2.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '156.84'
$ws.Range("E25").Value = '  +0.79%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.908'
$ws.Range("E26").Value = '  -1.68%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.52'
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.113'
$ws.Range("E28").Value = '  +5.20%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '117.69'
$ws.Range("E29").Value = '  +2.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.979'
$ws.Range("E30").Value = '  +0.77%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08938'
$ws.Range("E31").Value = '  +0.50%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.235'
$ws.Range("E32").Value = '  -2.64%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.256'
$ws.Range("E33").Value = '  +4.04%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7718'
$ws.Range("E34").Value = '  +2.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.665'
$ws.Range("E35").Value = '  +1.67%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.609'
$ws.Range("E36").Value = '  -3.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02052'
$ws.Range("E37").Value = '  +0.59%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.106'
$ws.Range("E38").Value = '  -1.28%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5516'
$ws.Range("E39").Value = '  -0.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05293'
$ws.Range("E40").Value = '  +0.65%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.992'
$ws.Range("E41").Value = '  +0.12%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.993'
$ws.Range("E42").Value = '  -0.36%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1527'
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.497'
$ws.Range("E44").Value = '  -1.20%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '109.97'
$ws.Range("E45").Value = '  +7.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.67'
$ws.Range("E46").Value = '  +0.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4834'
$ws.Range("E47").Value = '  -0.64%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.004'
$ws.Range("E48").Value = '  -0.70%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.649'
$ws.Range("E49").Value = '  -0.81%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '68.16'
$ws.Range("E50").Value = '  +1.32%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06078'
$ws.Range("E51").Value = '  -0.26%  '
